$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated figures from the 2020-08-18 Fonds de solidarite volet 2 data refresh.
# Values are stored as text in the source data (inlineStr), so a leading
# apostrophe forces Excel to keep numeric-looking strings as text rather
# than silently converting them to numbers.
$ws.Range("C2").Value = "'520"
$ws.Range("D2").Value = "'1241687.79"

$ws.Range("C4").Value = "'991"
$ws.Range("D4").Value = "'3472953.47"

$ws.Range("C6").Value = "'622"
$ws.Range("D6").Value = "'1968525.78"

$ws.Range("C8").Value = "'35"
$ws.Range("D8").Value = "'152144.45"

$ws.Range("C14").Value = "'215"
$ws.Range("D14").Value = "'578362.00"

$ws.Range("C16").Value = "'491"
$ws.Range("D16").Value = "'1791074.13"

$ws.Range("C28").Value = "'273"
$ws.Range("D28").Value = "'701657.45"

$ws.Range("C30").Value = "'545"
$ws.Range("D30").Value = "'2204872.70"

$ws.Range("C32").Value = "'382"
$ws.Range("D32").Value = "'1290613.57"

$ws.Range("C40").Value = "'141"
$ws.Range("D40").Value = "'395017.22"

$ws.Range("C41").Value = "'84"
$ws.Range("D41").Value = "'409909.98"

$ws.Range("C44").Value = "'7"
$ws.Range("D44").Value = "'38755.00"

$ws.Range("C45").Value = "'360"
$ws.Range("D45").Value = "'982867.74"

$ws.Range("C47").Value = "'592"
$ws.Range("D47").Value = "'2318908.99"

$ws.Range("C48").Value = "'401"
$ws.Range("D48").Value = "'1349177.16"

$ws.Range("C74").Value = "'383"
$ws.Range("D74").Value = "'967009.70"

$ws.Range("C76").Value = "'915"
$ws.Range("D76").Value = "'3160287.26"

$ws.Range("C77").Value = "'518"
$ws.Range("D77").Value = "'1703220.47"

$ws.Range("C92").Value = "'603"
$ws.Range("D92").Value = "'1469849.94"

$ws.Range("C94").Value = "'1087"
$ws.Range("D94").Value = "'3677253.80"

$ws.Range("C96").Value = "'1003"
$ws.Range("D96").Value = "'3060136.31"
